$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 553.61017
$ws.Range("J17").Value = 405.15686
$ws.Range("L17").Value = 1215.47058
$ws.Range("N17").Value = -1551.47058
$ws.Range("H41").Value = 318.55554
$ws.Range("J41").Value = 383.3846
$ws.Range("L41").Value = 383.3846
$ws.Range("N41").Value = -1263.3846
$ws.Range("H74").Value = 5630
$ws.Range("I74").Value = 4562.5
$ws.Range("K74").Value = 4562.5
$ws.Range("M74").Value = -3626.5
$ws.Range("H77").Value = 5630
$ws.Range("I77").Value = 4562.5
$ws.Range("K77").Value = 22812.5
$ws.Range("M77").Value = -18132.5
$ws.Range("H113").Value = 3704.7896
$ws.Range("I113").Value = 1472.5
$ws.Range("K113").Value = 1472.5
$ws.Range("M113").Value = 1781.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4159.375
$ws.Range("I74").Value = 4874.45
$ws.Range("J74").Value = 2967.5833
$ws.Range("K74").Value = 4874.45
$ws.Range("L74").Value = 2967.5833
$ws.Range("M74").Value = -4000.45
$ws.Range("N74").Value = -4715.5833
$ws.Range("H77").Value = 4159.375
$ws.Range("I77").Value = 4874.45
$ws.Range("J77").Value = 2967.5833
$ws.Range("K77").Value = 24372.25
$ws.Range("L77").Value = 14837.9165
$ws.Range("M77").Value = -20004.25
$ws.Range("N77").Value = -23573.9165
$ws.Range("H132").Value = 2157.3022
$ws.Range("I132").Value = 1454
$ws.Range("J132").Value = 2966.1
$ws.Range("K132").Value = 4362
$ws.Range("L132").Value = 8898.299999999999
$ws.Range("M132").Value = -1832
$ws.Range("N132").Value = -13958.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3983.6667
$ws.Range("I134").Value = 1287.8611
$ws.Range("J134").Value = 12071.083
$ws.Range("K134").Value = 3863.5833
$ws.Range("L134").Value = 36213.249
$ws.Range("M134").Value = -1328.5833
$ws.Range("N134").Value = -41283.249

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2415.6548
$ws.Range("I58").Value = 1559.4242
$ws.Range("J58").Value = 5555.1665
$ws.Range("K58").Value = 1559.4242
$ws.Range("L58").Value = 5555.1665
$ws.Range("M58").Value = -1356.4242
$ws.Range("N58").Value = -5961.1665
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H105").Value = 1663.4736
$ws.Range("J105").Value = 2248.1667
$ws.Range("L105").Value = 2248.1667
$ws.Range("N105").Value = -5742.1667
$ws.Range("H134").Value = 5084.8667
$ws.Range("I134").Value = 6326.5264
$ws.Range("J134").Value = 2940.182
$ws.Range("K134").Value = 18979.5792
$ws.Range("L134").Value = 8820.545999999998
$ws.Range("M134").Value = -16444.5792
$ws.Range("N134").Value = -13890.546
$ws.Range("H136").Value = 2415.6548
$ws.Range("I136").Value = 1559.4242
$ws.Range("J136").Value = 5555.1665
$ws.Range("K136").Value = 4678.2726
$ws.Range("L136").Value = 16665.4995
$ws.Range("M136").Value = -2128.2726
$ws.Range("N136").Value = -21765.4995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5008.3335
$ws.Range("I63").Value = 5250
$ws.Range("J63").Value = 4960
$ws.Range("K63").Value = 15750
$ws.Range("L63").Value = 14880
$ws.Range("M63").Value = -15001
$ws.Range("N63").Value = -16378
$ws.Range("H64").Value = 2522.5
$ws.Range("I64").Value = 1250
$ws.Range("J64").Value = 2946.6667
$ws.Range("K64").Value = 3750
$ws.Range("L64").Value = 8840.000100000001
$ws.Range("M64").Value = -3480
$ws.Range("N64").Value = -9380.000100000001
$ws.Range("H66").Value = 5008.3335
$ws.Range("I66").Value = 5250
$ws.Range("J66").Value = 4960
$ws.Range("K66").Value = 47250
$ws.Range("L66").Value = 44640
$ws.Range("M66").Value = -43506
$ws.Range("N66").Value = -52128
$ws.Range("H67").Value = 2522.5
$ws.Range("I67").Value = 1250
$ws.Range("J67").Value = 2946.6667
$ws.Range("K67").Value = 3750
$ws.Range("L67").Value = 8840.000100000001
$ws.Range("M67").Value = -2814
$ws.Range("N67").Value = -10712.0001
$ws.Range("H70").Value = 2837.3635
$ws.Range("I70").Value = 2504
$ws.Range("K70").Value = 7512
$ws.Range("M70").Value = -7197
$ws.Range("H73").Value = 2837.3635
$ws.Range("I73").Value = 2504
$ws.Range("K73").Value = 7512
$ws.Range("M73").Value = -6420
$ws.Range("H75").Value = 11499.5
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 11499.5
$ws.Range("K75").Value = 0
$ws.Range("M75").Value = 34498.5
$ws.Range("N75").Value = -36494.5
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 60000
$ws.Range("N76").Value = -60766
$ws.Range("H78").Value = 11499.5
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 11499.5
$ws.Range("K78").Value = 0
$ws.Range("M78").Value = 103495.5
$ws.Range("N78").Value = -113479.5
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 60000
$ws.Range("N79").Value = -62652
$ws.Range("H113").Value = 619.75
$ws.Range("J113").Value = 686
$ws.Range("L113").Value = 2058
$ws.Range("N113").Value = -6398
$ws.Range("H129").Value = 2307.8572
$ws.Range("I129").Value = 2259.8462
$ws.Range("J129").Value = 2385.875
$ws.Range("K129").Value = 6779.5386
$ws.Range("L129").Value = 7157.625
$ws.Range("M129").Value = -1779.5386
$ws.Range("N129").Value = -17157.625
$ws.Range("H131").Value = 774.54
$ws.Range("I131").Value = 371.66666
$ws.Range("J131").Value = 800.2553
$ws.Range("K131").Value = 1114.99998
$ws.Range("L131").Value = 2400.7659
$ws.Range("M131").Value = 3925.00002
$ws.Range("N131").Value = -12480.7659

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3168.3333
$ws.Range("I80").Value = 4005
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 4005
$ws.Range("L80").Value = 2750
$ws.Range("M80").Value = -3007
$ws.Range("N80").Value = -4746
$ws.Range("H83").Value = 3168.3333
$ws.Range("I83").Value = 4005
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 20025
$ws.Range("L83").Value = 13750
$ws.Range("M83").Value = -15033
$ws.Range("N83").Value = -23734
$ws.Range("H126").Value = 3322.36
$ws.Range("I126").Value = 2978.6494
$ws.Range("J126").Value = 4473.0435
$ws.Range("K126").Value = 8935.948199999999
$ws.Range("L126").Value = 13419.1305
$ws.Range("M126").Value = -6465.948199999999
$ws.Range("N126").Value = -18359.1305
$ws.Range("H132").Value = 2927.5
$ws.Range("I132").Value = 2133.75
$ws.Range("J132").Value = 3125.9375
$ws.Range("K132").Value = 6401.25
$ws.Range("L132").Value = 9377.8125
$ws.Range("M132").Value = -3871.25
$ws.Range("N132").Value = -14437.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1575.1052
$ws.Range("J46").Value = 1353.2858
$ws.Range("L46").Value = 1353.2858
$ws.Range("N46").Value = -1729.2858
$ws.Range("H61").Value = 1850.4445
$ws.Range("I61").Value = 1610.6666
$ws.Range("K61").Value = 1610.6666
$ws.Range("M61").Value = -1408.6666
$ws.Range("H68").Value = 642.4
$ws.Range("I68").Value = 642.4
$ws.Range("K68").Value = 642.4
$ws.Range("M68").Value = 106.6
$ws.Range("H71").Value = 642.4
$ws.Range("I71").Value = 642.4
$ws.Range("K71").Value = 3212
$ws.Range("M71").Value = 532
$ws.Range("H82").Value = 1948.68
$ws.Range("I82").Value = 661.4167
$ws.Range("J82").Value = 3136.923
$ws.Range("K82").Value = 661.4167
$ws.Range("L82").Value = 3136.923
$ws.Range("M82").Value = -300.4167
$ws.Range("N82").Value = -3858.923
$ws.Range("H85").Value = 1948.68
$ws.Range("I85").Value = 661.4167
$ws.Range("J85").Value = 3136.923
$ws.Range("K85").Value = 661.4167
$ws.Range("L85").Value = 3136.923
$ws.Range("M85").Value = 586.5833
$ws.Range("N85").Value = -5632.923
$ws.Range("H109").Value = 39750
$ws.Range("J109").Value = 39750
$ws.Range("L109").Value = 39750
$ws.Range("N109").Value = -42524
$ws.Range("H113").Value = 1850.4445
$ws.Range("I113").Value = 1610.6666
$ws.Range("K113").Value = 1610.6666
$ws.Range("M113").Value = 559.3334
$ws.Range("H132").Value = 22711.941
$ws.Range("I132").Value = 49580
$ws.Range("J132").Value = 11516.917
$ws.Range("K132").Value = 148740
$ws.Range("L132").Value = 34550.751
$ws.Range("M132").Value = -146210
$ws.Range("N132").Value = -39610.751
$ws.Range("H133").Value = 32526.2
$ws.Range("J133").Value = 32526.2
$ws.Range("L133").Value = 32526.2
$ws.Range("N133").Value = -37586.2
$ws.Range("H136").Value = 3643.8064
$ws.Range("I136").Value = 1720.7646
$ws.Range("J136").Value = 5978.9287
$ws.Range("K136").Value = 5162.293799999999
$ws.Range("L136").Value = 17936.7861
$ws.Range("M136").Value = -2612.293799999999
$ws.Range("N136").Value = -23036.7861

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2689.1428
$ws.Range("I132").Value = 1557.2941
$ws.Range("J132").Value = 7499.5
$ws.Range("K132").Value = 4671.8823
$ws.Range("L132").Value = 22498.5
$ws.Range("M132").Value = -2141.8823
$ws.Range("N132").Value = -27558.5
$ws.Range("H136").Value = 2451.587
$ws.Range("I136").Value = 1652.1177
$ws.Range("K136").Value = 4956.3531
$ws.Range("M136").Value = -2406.3531
